$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.426.23"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "2.092.93"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.5205"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "'0.4354"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").Value = "'54.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.36%  "

$ws.Range("D10").Value = "'0.08851"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "'1.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("D12").Value = "'24.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("D13").Value = "2.089.32"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").Value = "'6.667"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("D15").Value = "'7.656"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "'95.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").Value = "'0.06577"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").Value = "'19.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "'6.239"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").Value = "30.477.28"
$ws.Range("E23").Value = "  -1.14%  "

$ws.Range("D24").Value = "'12.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").Value = "'2.338"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.42%  "

$ws.Range("D26").Value = "2.333.02"
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").Value = "'22.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("D28").Value = "'2.558"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("D29").Value = "'162.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").Value = "'131.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("D31").Value = "'1.180"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").Value = "'0.1067"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "'1.650"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.87%  "

$ws.Range("D34").Value = "'6.118"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.90%  "

$ws.Range("D35").Value = "'3.887"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("D36").Value = "'9.984"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("D37").Value = "'0.02571"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("D38").Value = "'0.06797"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").Value = "'12.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").Value = "'5.434"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.46%  "

$ws.Range("D41").Value = "'0.2255"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'0.6876"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("D43").Value = "'1.263"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'0.6358"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").Value = "'13.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.82%  "

$ws.Range("D47").Value = "'2.191"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").Value = "'3.619"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("D49").Value = "'1.235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.58%  "

$ws.Range("D50").Value = "'1.238"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.93%  "

$ws.Range("D51").Value = "'81.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
